$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '48.166.44'
$ws.Range("E2").Value = '  +1.92%  '
$ws.Range("D3").Value = '2.515.18'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '321.65'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '109.32'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +1.38%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.530'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("E8").Value = '  +0.00%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.549'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +3.00%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '39.97'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +3.60%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '20.28'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +10.35%  '
$ws.Range("E12").Value = '  +1.26%  '
$ws.Range("E13").Value = '  +1.06%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '7.23'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = '2.910.82'
$ws.Range("E15").Value = '  +1.17%  '
$ws.Range("D16").Value = '2.518.48'
$ws.Range("E16").Value = '  +1.16%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.850'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").Value = '48.012.99'
$ws.Range("E18").Value = '  +1.82%  '
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '13.22'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +3.23%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.60'
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("E22").Value = '  -0.35%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '71.98'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +2.44%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '276.11'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +12.46%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.56'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.08%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '25.95'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +1.05%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.37'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +4.06%  '
$ws.Range("E29").Value = '  +1.38%  '
$ws.Range("E30").Value = '  +6.36%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '35.50'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.79%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '49.57'
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '19.69'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -4.18%  '
$ws.Range("E34").Value = '  +0.94%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +0.46%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.97'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("E38").Value = '  +1.00%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.97'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +1.78%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '122.88'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +3.76%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.112'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.51%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '2.22'
$c.Style = "Normal"
$ws.Range("E42").Value = '  -0.11%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '21.97'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("E44").Value = '  +1.62%  '
$ws.Range("D45").Value = '2.031.36'
$ws.Range("E45").Value = '  +2.43%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '3.12'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +3.33%  '
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.87'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +5.70%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.00'
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.45%  '
$ws.Range("E49").Value = '  -0.17%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '5.18'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +1.66%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '80.10'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +4.11%  '
